$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 307 (shifts existing rows 307..385 down to 308..386)
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Range("A307").Value = 7
$ws.Range("B307").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C307").Value = "Ñuble"
$ws.Range("D307").Value = 44798
$ws.Range("E307").Value = 16
$ws.Range("F307").Value = 100114001
$ws.Range("G307").Value = "Papa"
$ws.Range("H307").Value = "Patagonia"
$ws.Range("I307").Value = "1a (guarda)"
$ws.Range("J307").Value = 120
$ws.Range("K307").Value = 7000
$ws.Range("L307").Value = 7500
$ws.Range("M307").Value = 7250
$ws.Range("N307").Value = "`$/saco 25 kilos"
$ws.Range("O307").Value = "Provincia de Diguillín"
$ws.Range("P307").Value = 290
$ws.Range("Q307").Value = 25
$ws.Range("R307").Value = "Hortaliza"
